$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (column D) and Volume(1h) (column E) values per the commit diff
# (GitHub Actions "Updated cryptos list" run). Values are written as plain text
# so cells that look numeric (e.g. "244.93") stay text cells, matching the
# original inline-string layout - NumberFormat is forced to Text while writing
# and the cell style is reset to "Normal" afterwards so no formatting drifts.
$updates = @(
    @{ Row = 2; D = '43.222.35'; E = '  +4.60%  ' },
    @{ Row = 3; D = '2.251.04'; E = '  +3.34%  ' },
    @{ Row = 4; D = $null; E = '  +0.07%  ' },
    @{ Row = 5; D = '244.93'; E = '  +3.06%  ' },
    @{ Row = 6; D = '0.617'; E = '  +1.35%  ' },
    @{ Row = 7; D = '76.05'; E = '  +8.52%  ' },
    @{ Row = 9; D = $null; E = '  +6.42%  ' },
    @{ Row = 10; D = '41.37'; E = '  +4.71%  ' },
    @{ Row = 11; D = $null; E = '  +1.83%  ' },
    @{ Row = 12; D = $null; E = '  +4.51%  ' },
    @{ Row = 13; D = '0.101'; E = '  +0.51%  ' },
    @{ Row = 14; D = '2.587.92'; E = '  +3.45%  ' },
    @{ Row = 15; D = '14.67'; E = '  +2.61%  ' },
    @{ Row = 16; D = '2.245.58'; E = '  +3.99%  ' },
    @{ Row = 17; D = '0.804'; E = '  +1.60%  ' },
    @{ Row = 18; D = '43.114.11'; E = '  +4.90%  ' },
    @{ Row = 19; D = $null; E = '  +5.40%  ' },
    @{ Row = 20; D = '71.49'; E = '  +1.27%  ' },
    @{ Row = 21; D = $null; E = '  +2.10%  ' },
    @{ Row = 22; D = '10.02'; E = '  +5.89%  ' },
    @{ Row = 23; D = '230.79'; E = '  +2.18%  ' },
    @{ Row = 24; D = $null; E = '  +14.94%  ' },
    @{ Row = 25; D = $null; E = '  +0.03%  ' },
    @{ Row = 26; D = '10.97'; E = '  +1.89%  ' },
    @{ Row = 27; D = '3.51'; E = '  +0.96%  ' },
    @{ Row = 28; D = '39.40'; E = '  +30.44%  ' },
    @{ Row = 29; D = $null; E = '  +2.42%  ' },
    @{ Row = 30; D = $null; E = '  +2.14%  ' },
    @{ Row = 31; D = '173.54'; E = '  +3.53%  ' },
    @{ Row = 32; D = '20.39'; E = '  +2.31%  ' },
    @{ Row = 33; D = $null; E = '  +4.66%  ' },
    @{ Row = 34; D = '5.36'; E = '  +4.54%  ' },
    @{ Row = 35; D = $null; E = '  +1.60%  ' },
    @{ Row = 36; D = '0.111'; E = '  +9.03%  ' },
    @{ Row = 37; D = $null; E = '  +6.41%  ' },
    @{ Row = 38; D = '0.0336'; E = '  +19.04%  ' },
    @{ Row = 39; D = '13.29'; E = '  +12.09%  ' },
    @{ Row = 40; D = $null; E = '  +3.40%  ' },
    @{ Row = 41; D = '5.54'; E = '  +3.07%  ' },
    @{ Row = 42; D = '0.204'; E = '  +7.27%  ' },
    @{ Row = 43; D = '60.06'; E = '  +1.87%  ' },
    @{ Row = 44; D = '105.97'; E = '  +8.86%  ' },
    @{ Row = 45; D = '8.75'; E = '  +5.50%  ' },
    @{ Row = 46; D = '0.482'; E = '  +29.75%  ' },
    @{ Row = 47; D = '0.0996'; E = '  +2.94%  ' },
    @{ Row = 48; D = '2.41'; E = '  +10.13%  ' },
    @{ Row = 49; D = $null; E = '  +3.22%  ' },
    @{ Row = 50; D = '1.15'; E = '  +2.39%  ' },
    @{ Row = 51; D = '2.462.83'; E = '  +3.53%  ' }
)

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.D) {
        Set-TextValue $ws.Range("D$row") $u.D
    }
    if ($null -ne $u.E) {
        Set-TextValue $ws.Range("E$row") $u.E
    }
}
